$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Players should still have gravity after being frozen (with exceptions), so the two
# "ground the player" TODO entries are no longer needed - remove their rows entirely.
$target1 = $ws.Range("A1:A200").Find("somehow ground player after teleporting")
if ($target1 -ne $null) {
    $target1.EntireRow.Delete()
}

$target2 = $ws.Range("A1:A200").Find("somehow smoothly ground player during freeze")
if ($target2 -ne $null) {
    $target2.EntireRow.Delete()
}

$wb.Save()
